# ScenarioTemplate / UIDMatrix sheet update
# - Add a handful of per-row "UID" helper cells (col C/I/O) that mirror the
#   style already used in column A (s=25, red text on yellow fill).
# - Expand the "Harvested wood products" (HWP) block (rows 61-65) into a
#   more detailed eleven-row block (rows 61-71) that breaks Sawnwood, Wood
#   panels and Paper & paperboard each into domestic/exported/total lines,
#   each carrying its own UID cell in column C.
# - Refresh the sheet view (scroll position / selection) and add an explicit
#   page setup so printing uses A4-ish "paperSize 9 / portrait" settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New UID cells in the top matrix (rows 9-15), styled like column A
#    (copy format only from a cell that already carries style s=25).
# ---------------------------------------------------------------------
function Copy-StyleOnly($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$styleSrc = $ws.Range("A9")

Copy-StyleOnly $styleSrc $ws.Range("I9")
$ws.Range("I9").Value2 = "176984AA-39DD-46BD-8783-2632BEF3C520"

Copy-StyleOnly $styleSrc $ws.Range("C10")
$ws.Range("C10").Value2 = "E1B2A0CF-5F9F-445E-9715-C274A0CD4A26"

Copy-StyleOnly $styleSrc $ws.Range("I10")
$ws.Range("I10").Value2 = "0A0CAA48-DB6F-412A-AFBD-8F078B1AF8A6"

Copy-StyleOnly $styleSrc $ws.Range("O10")
$ws.Range("O10").Value2 = "FFD8E79B-1DA0-4399-910A-6E875F1A8F58"

Copy-StyleOnly $styleSrc $ws.Range("C14")
$ws.Range("C14").Value2 = "33B17FCF-CEF0-4C52-A083-41E77975CC17"

Copy-StyleOnly $styleSrc $ws.Range("C15")
$ws.Range("C15").Value2 = "1632C1F2-832E-48D5-BA76-AA1DFAA643DC"

# ---------------------------------------------------------------------
# 2) Expand the HWP block: rows 61-65 (5 rows) become rows 61-71 (11 rows).
#    Insert 4 new blank rows at 66 - the 3 already-blank separator rows
#    that used to sit right below the block (old rows 66-68) absorb the
#    rest of the growth, so everything from the old row 69 onward ends up
#    shifted down by exactly 4 rows (old 69 -> new 73, old 126 -> new 130).
# ---------------------------------------------------------------------
$ws.Rows("66:69").Insert()

# Style template for the detail rows/cells (copy from existing A62 / B62,
# which already carry cell style s=31 inside the s=13 custom-formatted
# HWP rows).
$cellStyleA = $ws.Range("A62")
$cellStyleB = $ws.Range("B62")
$cellStyleC = $ws.Range("C63")

function Set-HwpCell($addr, $template, $value) {
    Copy-StyleOnly $template $ws.Range($addr)
    if ($null -ne $value) {
        $ws.Range($addr).Value2 = $value
    }
}

# Row 61: section header "Harvested wood products" (A), blank styled B/C
Set-HwpCell "A61" $cellStyleA $null
$ws.Range("A61").Value2 = "Harvested wood products"
Set-HwpCell "B61" $cellStyleB $null
Set-HwpCell "C61" $cellStyleC $null

# Row 62: Sawnwood domestic
Set-HwpCell "A62" $cellStyleA 27
Set-HwpCell "B62" $cellStyleB "Sawnwood domestic"
Set-HwpCell "C62" $cellStyleC "0B29C52E-CF04-46A9-AF23-BA996645E547"

# Row 63: Sawnwood exported
Set-HwpCell "A63" $cellStyleA 28
Set-HwpCell "B63" $cellStyleB "Sawnwood exported"
Set-HwpCell "C63" $cellStyleC "054E0C6F-5E84-4411-9022-B04049A2C6D0"

# Row 64: Sawnwood total
Set-HwpCell "A64" $cellStyleA $null
Set-HwpCell "B64" $cellStyleB "Sawnwood total"
Set-HwpCell "C64" $cellStyleC $null

# Row 65: Wood panels domestic
Set-HwpCell "A65" $cellStyleA 29
Set-HwpCell "B65" $cellStyleB "Wood panels domestic"
Set-HwpCell "C65" $cellStyleC "CA184B1C-EB7C-42CB-897C-57AC319EAF39"

# Row 66: Wood panels exported
Set-HwpCell "A66" $cellStyleA 30
Set-HwpCell "B66" $cellStyleB "Wood panels exported"
Set-HwpCell "C66" $cellStyleC "42476E95-C4C0-4787-8E7B-2CDA36458710"

# Row 67: Wood panels total
Set-HwpCell "A67" $cellStyleA $null
Set-HwpCell "B67" $cellStyleB "Wood panels total"
Set-HwpCell "C67" $cellStyleC $null

# Row 68: Paper and paperboard domestic
Set-HwpCell "A68" $cellStyleA 31
Set-HwpCell "B68" $cellStyleB "Paper and paperboard domestic"
Set-HwpCell "C68" $cellStyleC "185C43A5-D840-460E-8163-8B2F7BBB952C"

# Row 69: Paper and paperboard exported
Set-HwpCell "A69" $cellStyleA 32
Set-HwpCell "B69" $cellStyleB "Paper and paperboard exported"
Set-HwpCell "C69" $cellStyleC "3BBB8E1F-6E69-4551-A0F4-527BFBBBDBFB"

# Row 70: Paper and paperboard total
Set-HwpCell "A70" $cellStyleA $null
Set-HwpCell "B70" $cellStyleB "Paper and paperboard total"
Set-HwpCell "C70" $cellStyleC $null

# Row 71: HWP total
Set-HwpCell "A71" $cellStyleA $null
Set-HwpCell "B71" $cellStyleB "HWP total"
Set-HwpCell "C71" $cellStyleC $null

# ---------------------------------------------------------------------
# 3) View / selection / print setup
# ---------------------------------------------------------------------
$ws.Range("B81").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1

$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
